# Apply updated "想去人数" (want-to-go count) figures and a corrected venue
# name for the 苏州-漫展信息 workbook. The same data is duplicated across the
# "展览" and "全部类型" worksheets, so the same edits are applied to both.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 11451
    $ws.Range("F3").Value = 10868
    $ws.Range("F4").Value = 602
    $ws.Range("F8").Value = 54
    $ws.Range("F9").Value = 37
    $ws.Range("F11").Value = 10591
    $ws.Range("F12").Value = 4088
    $ws.Range("F14").Value = 2453
    $ws.Range("D15").Value = "金山南路影视城 木渎影视城会展中心"
    $ws.Range("F17").Value = 107
    $ws.Range("F18").Value = 418
    $ws.Range("F19").Value = 11096
    $ws.Range("F20").Value = 10854
}
